{"js": "// Replace the body of the \"problem statement\" paragraph with the\n// rewritten / shortened text, keeping the same run formatting.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldLead = \"O problema \u00e9 a dificuldade da controlar os hor\u00e1rios\";\nconst newText =\n  \"O problema da dificuldade de controlar os hor\u00e1rios e servi\u00e7os marcados, \" +\n  \"afeta os funcion\u00e1rios, devido perca de credibilidade com os clientes, \" +\n  \"insatisfa\u00e7\u00e3o por demora no atendimento e reclama\u00e7\u00f5es. \";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(oldLead) === 0) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the problem-statement paragraph.\");\n}\n\ntarget.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"problem statement\" paragraph by searching for its\n# distinctive opening words, then rewrite/shorten its text in place\n# (the run formatting -- Arial, black, 12pt -- is preserved because we\n# only change the Range.Text, not the run properties).\n$searchRange = $d.Content\n$searchRange.Find.MatchCase = $true\n$searchRange.Find.Text = \"O problema \u00e9 a dificuldade\"\n$found = $searchRange.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not locate the problem-statement paragraph.\"\n}\n\n$targetParagraph = $searchRange.Paragraphs(1)\n$newText = \"O problema da dificuldade de controlar os hor\u00e1rios e servi\u00e7os marcados, afeta os funcion\u00e1rios, devido perca de credibilidade com os clientes, insatisfa\u00e7\u00e3o por demora no atendimento e reclama\u00e7\u00f5es. \"\n\n$targetParagraph.Range.Text = $newText\n"}
